$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 49 and 50: BabyDogeCoin / PancakeSwap swapped order
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D49") "0.00000000335"
Set-TextValue $ws.Range("E49") "  -1.73%  "

$ws.Range("B50").Value = "PancakeSwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D50") "3.389"
Set-TextValue $ws.Range("E50") "  -3.12%  "

# Remaining rows: D/E price + volume updates
Set-TextValue $ws.Range("D2") "28.424.60"
Set-TextValue $ws.Range("E2") "  -3.10%  "
Set-TextValue $ws.Range("D3") "1.954.48"
Set-TextValue $ws.Range("E3") "  -1.34%  "
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("E4") "  -1.30%  "
Set-TextValue $ws.Range("D5") "320.04"
Set-TextValue $ws.Range("E5") "  -2.59%  "
Set-TextValue $ws.Range("E6") "  -1.08%  "
Set-TextValue $ws.Range("D7") "0.4764"
Set-TextValue $ws.Range("E7") "  -4.92%  "
Set-TextValue $ws.Range("D8") "0.4039"
Set-TextValue $ws.Range("E8") "  -4.09%  "
Set-TextValue $ws.Range("D9") "53.43"
Set-TextValue $ws.Range("E9") "  -0.65%  "
Set-TextValue $ws.Range("D10") "0.08412"
Set-TextValue $ws.Range("E10") "  -5.18%  "
Set-TextValue $ws.Range("D11") "1.055"
Set-TextValue $ws.Range("E11") "  -4.34%  "
Set-TextValue $ws.Range("D12") "22.31"
Set-TextValue $ws.Range("E12") "  -3.14%  "
Set-TextValue $ws.Range("D13") "1.944.60"
Set-TextValue $ws.Range("E13") "  -9.55%  "
Set-TextValue $ws.Range("D14") "7.576"
Set-TextValue $ws.Range("E14") "  -3.84%  "
Set-TextValue $ws.Range("D15") "6.145"
Set-TextValue $ws.Range("E15") "  -4.30%  "
Set-TextValue $ws.Range("D16") "1.008"
Set-TextValue $ws.Range("E16") "  -0.99%  "
Set-TextValue $ws.Range("E17") "  -2.87%  "
Set-TextValue $ws.Range("D18") "0.00001069"
Set-TextValue $ws.Range("E18") "  -3.07%  "
Set-TextValue $ws.Range("D19") "0.06578"
Set-TextValue $ws.Range("E19") "  -2.70%  "
Set-TextValue $ws.Range("D20") "18.51"
Set-TextValue $ws.Range("E20") "  -3.99%  "
Set-TextValue $ws.Range("D21") "1.005"
Set-TextValue $ws.Range("E21") "  -1.06%  "
Set-TextValue $ws.Range("D22") "5.813"
Set-TextValue $ws.Range("E22") "  -2.10%  "
Set-TextValue $ws.Range("D23") "28.435.93"
Set-TextValue $ws.Range("E23") "  -3.42%  "
Set-TextValue $ws.Range("D24") "11.47"
Set-TextValue $ws.Range("E24") "  -3.72%  "
Set-TextValue $ws.Range("D25") "2.287"
Set-TextValue $ws.Range("E25") "  -1.20%  "
Set-TextValue $ws.Range("D26") "2.167.75"
Set-TextValue $ws.Range("E26") "  -7.76%  "
Set-TextValue $ws.Range("D27") "155.03"
Set-TextValue $ws.Range("E27") "  -0.77%  "
Set-TextValue $ws.Range("D28") "20.17"
Set-TextValue $ws.Range("E28") "  -2.17%  "
Set-TextValue $ws.Range("D29") "5.921"
Set-TextValue $ws.Range("E29") "  -4.32%  "
Set-TextValue $ws.Range("D30") "2.155"
Set-TextValue $ws.Range("E30") "  -5.71%  "
Set-TextValue $ws.Range("D31") "123.56"
Set-TextValue $ws.Range("E31") "  -2.34%  "
Set-TextValue $ws.Range("D32") "0.9797"
Set-TextValue $ws.Range("E32") "  -6.51%  "
Set-TextValue $ws.Range("D33") "0.09606"
Set-TextValue $ws.Range("E33") "  -2.86%  "
Set-TextValue $ws.Range("D34") "1.451"
Set-TextValue $ws.Range("E34") "  -4.04%  "
Set-TextValue $ws.Range("D35") "5.599"
Set-TextValue $ws.Range("E35") "  -3.22%  "
Set-TextValue $ws.Range("E36") "  -3.32%  "
Set-TextValue $ws.Range("D37") "8.897"
Set-TextValue $ws.Range("E37") "  -2.93%  "
Set-TextValue $ws.Range("D38") "0.02327"
Set-TextValue $ws.Range("E38") "  -4.33%  "
Set-TextValue $ws.Range("D39") "0.06222"
Set-TextValue $ws.Range("E39") "  -1.99%  "
Set-TextValue $ws.Range("D40") "1.247"
Set-TextValue $ws.Range("E40") "  -3.17%  "
Set-TextValue $ws.Range("D41") "0.6197"
Set-TextValue $ws.Range("E41") "  -4.51%  "
Set-TextValue $ws.Range("D42") "11.12"
Set-TextValue $ws.Range("E42") "  -3.48%  "
Set-TextValue $ws.Range("E43") "  -0.96%  "
Set-TextValue $ws.Range("D44") "0.1920"
Set-TextValue $ws.Range("E44") "  -5.10%  "
Set-TextValue $ws.Range("D45") "1.351"
Set-TextValue $ws.Range("E45") "  +4.65%  "
Set-TextValue $ws.Range("D46") "0.5948"
Set-TextValue $ws.Range("E46") "  -4.90%  "
Set-TextValue $ws.Range("D47") "13.03"
Set-TextValue $ws.Range("E47") "  -3.07%  "
Set-TextValue $ws.Range("D48") "2.058"
Set-TextValue $ws.Range("E48") "  -6.15%  "
Set-TextValue $ws.Range("D51") "0.06804"
Set-TextValue $ws.Range("E51") "  -1.68%  "
